$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily COVID overview rows for the UK, 2021-03-05 through 2021-03-21
# (date, areaType, areaCode, areaName, cumCases, newCases, newDeaths28Days, cumDeaths28Days)
$data = @(
    @("2021-03-05", "overview", "K02000001", "United Kingdom", 4207304, 5947, 236, 124261),
    @("2021-03-06", "overview", "K02000001", "United Kingdom", 4213343, 6040, 158, 124419),
    @("2021-03-07", "overview", "K02000001", "United Kingdom", 4218520, 5177, 82, 124501),
    @("2021-03-08", "overview", "K02000001", "United Kingdom", 4223232, 4712, 65, 124566),
    @("2021-03-09", "overview", "K02000001", "United Kingdom", 4228998, 5766, 231, 124797),
    @("2021-03-10", "overview", "K02000001", "United Kingdom", 4234924, 5926, 190, 124987),
    @("2021-03-11", "overview", "K02000001", "United Kingdom", 4241677, 6753, 181, 125168),
    @("2021-03-12", "overview", "K02000001", "United Kingdom", 4248286, 6609, 175, 125343),
    @("2021-03-13", "overview", "K02000001", "United Kingdom", 4253820, 5534, 121, 125464),
    @("2021-03-14", "overview", "K02000001", "United Kingdom", 4258438, 4618, 52, 125516),
    @("2021-03-15", "overview", "K02000001", "United Kingdom", 4263527, 5089, 64, 125580),
    @("2021-03-16", "overview", "K02000001", "United Kingdom", 4268821, 5294, 110, 125690),
    @("2021-03-17", "overview", "K02000001", "United Kingdom", 4274579, 5758, 141, 125831),
    @("2021-03-18", "overview", "K02000001", "United Kingdom", 4280882, 6303, 95, 125926),
    @("2021-03-19", "overview", "K02000001", "United Kingdom", 4285684, 4802, 101, 126026),
    @("2021-03-20", "overview", "K02000001", "United Kingdom", 4291271, 5587, 96, 126122),
    @("2021-03-21", "overview", "K02000001", "United Kingdom", 4296583, 5312, 33, 126155),
)

$startRow = 206
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    # Leading apostrophe forces the date to be stored as text, matching the
    # existing rows (column A holds plain "yyyy-mm-dd" strings, not date serials)
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Restore the view: scroll so row 186 is at the top and select A208,
# mirroring where the sheet was left after appending the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 186
$win.ScrollColumn = 1
$ws.Range("A208").Select()
